$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task 2")

# New "Status" column (D) with a Pass/fail marker next to each test case.
# Header cell — bold header style, centered.
$ws.Range("D1").Value = "Status"
$ws.Range("D1").HorizontalAlignment = -4108

# Row 2 (first test case) — "Pass", centered.
$ws.Range("D2").Value = "Pass"
$ws.Range("D2").HorizontalAlignment = -4108

# Row 3 (spacer row) — centered, no value.
$ws.Range("D3").HorizontalAlignment = -4108

# Row 4 (second test case) — "Pass", centered + wrap text (row is tall/wrapped).
$ws.Range("D4").Value = "Pass"
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").WrapText = $true

# Row 5 (spacer row) — quote-prefixed empty text, centered.
$ws.Range("D5").Value = "'"
$ws.Range("D5").HorizontalAlignment = -4108

# Row 6 (third test case) — "Pass", centered.
$ws.Range("D6").Value = "Pass"
$ws.Range("D6").HorizontalAlignment = -4108

# Row 7 (spacer row) — centered, no value.
$ws.Range("D7").HorizontalAlignment = -4108

# Select the whole new column, matching the saved selection state.
$ws.Columns("D:D").Select()
